$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update title text (A1, merged A1:G2) ---
$ws.Range("A1").Value = "Table 1. Oscar Elton Sette (SE) shipboard survey information for data used in this study with geographic boxes used to define frontal areas. All transects were from 26°-36°N along 158°W. Frontal positions were determined using CTD profile data."

# --- Clear old Latitude/Longitude/Transect length/Frontal Positions columns (C:G) for header rows ---
$ws.Range("C3:G7").Clear()

# --- Rebuild header row 3: Cruise / Dates / Frontal Positions (merged C3:D3) ---
$ws.Range("A3").Value = "Cruise "
$ws.Range("B3").Value = "Dates"
$ws.Range("C3").Value = "Frontal Positions"
$ws.Range("C3:D3").Merge()

# --- Row 4 sub-headers: STF / TZCF ---
$ws.Range("C4").Value = "STF"
$ws.Range("D4").Value = "TZCF"

# --- Data rows ---
$ws.Range("A5").Value = "SE-08-02"
$ws.Range("B5").Value = "26 March-3 April 2008"
$ws.Range("C5").Value = "32°15'-32°45'N"
$ws.Range("D5").Value = "34°15'-35°45'N"

$ws.Range("A6").Value = "SE-09-02"
$ws.Range("B6").Value = "18 March-23 March 2009"
$ws.Range("C6").Value = "31°15'-32°15'N"
$ws.Range("D6").Value = "35°00'-36°00'N"

$ws.Range("A7").Value = "SE-11-02"
$ws.Range("B7").Value = "10 March-23 March 2011"
$ws.Range("C7").Value = "31°15'-32°15'N"
$ws.Range("D7").Value = "33°15'-34°15'N"

# --- Column D width change ---
$ws.Columns.Item(4).ColumnWidth = 12

# --- Selection change ---
$ws.Range("A8").Select()
